# Auto-generated script applying scheduled-runner value updates
# to the Typhon_Profits Leve profit-tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 316.84616
$ws.Range("I53").Value = 306.16666
$ws.Range("K53").Value = 306.16666
$ws.Range("M53").Value = 330.83334
$ws.Range("H64").Value = 3163.158
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 3566.6667
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 3566.6667
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -4062.6667
$ws.Range("H67").Value = 3163.158
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 3566.6667
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 3566.6667
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -5282.6667
$ws.Range("H137").Value = 1737.4166
$ws.Range("I137").Value = 1271.2858
$ws.Range("J137").Value = 2390
$ws.Range("K137").Value = 3813.8574
$ws.Range("L137").Value = 7170
$ws.Range("M137").Value = -1263.8574
$ws.Range("N137").Value = -12270

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5136.857
$ws.Range("I61").Value = 6523.5557
$ws.Range("J61").Value = 4096.8335
$ws.Range("K61").Value = 6523.5557
$ws.Range("L61").Value = 4096.8335
$ws.Range("M61").Value = -6311.5557
$ws.Range("N61").Value = -4520.8335
$ws.Range("H97").Value = 1602.1818
$ws.Range("I97").Value = 1602.1818
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1602.1818
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1106.1818
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 2098.75
$ws.Range("I110").Value = 2098.75
$ws.Range("K110").Value = 2098.75
$ws.Range("M110").Value = -53.75
$ws.Range("H136").Value = 5136.857
$ws.Range("I136").Value = 6523.5557
$ws.Range("J136").Value = 4096.8335
$ws.Range("K136").Value = 19570.6671
$ws.Range("L136").Value = 12290.5005
$ws.Range("M136").Value = -17020.6671
$ws.Range("N136").Value = -17390.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 58824344
$ws.Range("I64").Value = 166668270
$ws.Range("J64").Value = 381.0909
$ws.Range("K64").Value = 166668270
$ws.Range("L64").Value = 381.0909
$ws.Range("M64").Value = -166668045
$ws.Range("N64").Value = -831.0908999999999
$ws.Range("H67").Value = 58824344
$ws.Range("I67").Value = 166668270
$ws.Range("J67").Value = 381.0909
$ws.Range("K67").Value = 166668270
$ws.Range("L67").Value = 381.0909
$ws.Range("M67").Value = -166667490
$ws.Range("N67").Value = -1941.0909
$ws.Range("H94").Value = 4124.9
$ws.Range("I94").Value = 1562.25
$ws.Range("J94").Value = 5833.3335
$ws.Range("K94").Value = 1562.25
$ws.Range("L94").Value = 5833.3335
$ws.Range("M94").Value = -1111.25
$ws.Range("N94").Value = -6735.3335
$ws.Range("H97").Value = 8088.3335
$ws.Range("I97").Value = 3642.1667
$ws.Range("J97").Value = 16980.666
$ws.Range("K97").Value = 3642.1667
$ws.Range("L97").Value = 16980.666
$ws.Range("M97").Value = -2651.1667
$ws.Range("N97").Value = -18962.666
$ws.Range("H107").Value = 1159.2727
$ws.Range("I107").Value = 1139.1111
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 1139.1111
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 780.8888999999999
$ws.Range("N107").Value = -5090
$ws.Range("H134").Value = 4006.7
$ws.Range("I134").Value = 4164.9473
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 12494.8419
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -9959.841899999999
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10857.25
$ws.Range("I31").Value = 13443
$ws.Range("J31").Value = 3100
$ws.Range("K31").Value = 13443
$ws.Range("L31").Value = 3100
$ws.Range("M31").Value = -13148
$ws.Range("N31").Value = -3690
$ws.Range("H34").Value = 10857.25
$ws.Range("I34").Value = 13443
$ws.Range("J34").Value = 3100
$ws.Range("K34").Value = 13443
$ws.Range("L34").Value = 3100
$ws.Range("M34").Value = -13241
$ws.Range("N34").Value = -3504
$ws.Range("H135").Value = 48970
$ws.Range("J135").Value = 48970
$ws.Range("L135").Value = 48970
$ws.Range("N135").Value = -59110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 899
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 899
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 4167.037
$ws.Range("I107").Value = 5164.1904
$ws.Range("J107").Value = 677
$ws.Range("K107").Value = 15492.5712
$ws.Range("L107").Value = 2031
$ws.Range("M107").Value = -13572.5712
$ws.Range("N107").Value = -5871
$ws.Range("H131").Value = 808.46
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 824.18555
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2472.55665
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12552.55665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2598.0667
$ws.Range("I97").Value = 1514.1428
$ws.Range("J97").Value = 3546.5
$ws.Range("K97").Value = 1514.1428
$ws.Range("L97").Value = 3546.5
$ws.Range("M97").Value = -1018.1428
$ws.Range("N97").Value = -4538.5
$ws.Range("H126").Value = 4291.838
$ws.Range("I126").Value = 3445.75
$ws.Range("J126").Value = 5853.846
$ws.Range("K126").Value = 10337.25
$ws.Range("L126").Value = 17561.538
$ws.Range("M126").Value = -7867.25
$ws.Range("N126").Value = -22501.538
$ws.Range("H132").Value = 34577.06
$ws.Range("I132").Value = 6421.2
$ws.Range("J132").Value = 74799.71000000001
$ws.Range("K132").Value = 19263.6
$ws.Range("L132").Value = 224399.13
$ws.Range("M132").Value = -16733.6
$ws.Range("N132").Value = -229459.13
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6362.5
$ws.Range("I7").Value = 6366.6665
$ws.Range("J7").Value = 6350
$ws.Range("K7").Value = 6366.6665
$ws.Range("L7").Value = 6350
$ws.Range("M7").Value = -6254.6665
$ws.Range("N7").Value = -6574
$ws.Range("H68").Value = 4502.8184
$ws.Range("I68").Value = 3363.3333
$ws.Range("J68").Value = 5870.2
$ws.Range("K68").Value = 3363.3333
$ws.Range("L68").Value = 5870.2
$ws.Range("M68").Value = -2614.3333
$ws.Range("N68").Value = -7368.2
$ws.Range("H71").Value = 4502.8184
$ws.Range("I71").Value = 3363.3333
$ws.Range("J71").Value = 5870.2
$ws.Range("K71").Value = 16816.6665
$ws.Range("L71").Value = 29351
$ws.Range("M71").Value = -13072.6665
$ws.Range("N71").Value = -36839
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("H126").Value = 6362.5
$ws.Range("I126").Value = 6366.6665
$ws.Range("J126").Value = 6350
$ws.Range("K126").Value = 19099.9995
$ws.Range("L126").Value = 19050
$ws.Range("M126").Value = -16629.9995
$ws.Range("N126").Value = -23990
$ws.Range("H132").Value = 432107.66
$ws.Range("I132").Value = 549078
$ws.Range("J132").Value = 3216.3333
$ws.Range("K132").Value = 1647234
$ws.Range("L132").Value = 9648.999899999999
$ws.Range("M132").Value = -1644704
$ws.Range("N132").Value = -14708.9999

